$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns (B:D) in front of the existing (empty) B:D,
# pushing nothing in column A around - only adds blank columns.
$ws.Range("B1:D1").EntireColumn.Insert()

# Row 1 - familyName (texto)
$ws.Range("A1").Value = "familyName"
$ws.Range("C1").Value = " "
$ws.Range("D1").Value = "texto"

# Row 2 - givenName (texto)
$ws.Range("A2").Value = "givenName"
$ws.Range("C2").Value = " "
$ws.Range("D2").Value = "texto"

# Row 3 - theDate (fecha)
$ws.Range("A3").Value = "theDate"
$ws.Range("B3").Value = "1/1/0001 00:00:00"
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = "fecha"

# Row 4 - TextField1 (texto)
$ws.Range("A4").Value = "TextField1"
$ws.Range("C4").Value = " "
$ws.Range("D4").Value = "texto"

# Row 5 - DateTimeField1 (fecha)
$ws.Range("A5").Value = "DateTimeField1"
$ws.Range("B5").Value = "1/1/0001 00:00:00"
$ws.Range("C5").Value = " "
$ws.Range("D5").Value = "fecha"
